$wb = $excel.ActiveWorkbook

# Update "Ready for handoff" -> "In Translation" everywhere it appears
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Autofit affected columns
$overview.Range("E:F").Columns.AutoFit() | Out-Null
$zhcn.Range("C:C").Columns.AutoFit() | Out-Null
$dede.Range("C:C").Columns.AutoFit() | Out-Null
